$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 'Depth'
$ws.Cells.Item(1, 3).Value = 'DRM'
$ws.Cells.Item(1, 4).Value = 'Month'
$ws.Cells.Item(1, 5).Value = 'Depth:DRM'
$ws.Cells.Item(1, 6).Value = 'Depth:Month'
$ws.Cells.Item(1, 7).Value = 'DRM:Month'
$ws.Cells.Item(2, 2).Value = [double]"-0.05869947718450887"
$ws.Cells.Item(2, 3).Value = [double]"-0.143497918224461"
$ws.Cells.Item(2, 4).Value = '+'
$ws.Cells.Item(2, 5).Value = [double]"0.06179045099957722"
$ws.Cells.Item(2, 7).Value = '+'
$ws.Cells.Item(2, 8).Value = [double]"0.7929473038810541"
$ws.Cells.Item(2, 13).Value = [double]"0.7662212737020859"
$ws.Cells.Item(2, 14).Value = [double]"0.7451659124689896"
$ws.Cells.Item(3, 1).Value = [double]"4.033399788238579"
$ws.Cells.Item(3, 2).Value = [double]"-0.08077016576720128"
$ws.Cells.Item(3, 3).Value = [double]"-0.1812764453405394"
$ws.Cells.Item(3, 4).Value = '+'
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(3, 7).Value = '+'
$ws.Cells.Item(4, 2).Value = [double]"0.01527448987367492"
$ws.Cells.Item(4, 3).Value = [double]"-0.09680689775510189"
$ws.Cells.Item(4, 4).Value = '+'
$ws.Cells.Item(4, 5).Value = [double]"0.08914883112156308"
$ws.Cells.Item(4, 6).ClearContents()
$ws.Cells.Item(4, 7).Value = '+'
$ws.Cells.Item(4, 8).Value = [double]"0.7193981157696526"
$ws.Cells.Item(4, 13).Value = [double]"0.02954966487824302"
$ws.Cells.Item(4, 14).Value = [double]"0.6674348038751439"
$ws.Cells.Item(5, 2).ClearContents()
$ws.Cells.Item(5, 3).Value = [double]"-0.1435795843332341"
$ws.Cells.Item(5, 4).Value = '+'
$ws.Cells.Item(5, 6).ClearContents()
$ws.Cells.Item(5, 7).Value = '+'
$ws.Cells.Item(5, 8).Value = [double]"0.6152040256792846"
$ws.Cells.Item(5, 12).Value = [double]"10.67363976377604"
$ws.Cells.Item(5, 13).Value = [double]"0.003686402739443487"
$ws.Cells.Item(5, 14).Value = [double]"0.5753975455771416"
$ws.Cells.Item(6, 2).Value = [double]"-0.01021751290752975"
$ws.Cells.Item(6, 3).Value = [double]"-0.03622644666262292"
$ws.Cells.Item(6, 4).Value = '+'
$ws.Cells.Item(6, 5).Value = [double]"0.08086217946940651"
$ws.Cells.Item(6, 6).Value = '+'
$ws.Cells.Item(6, 7).ClearContents()
$ws.Cells.Item(6, 10).Value = [double]"28.25201495402969"
$ws.Cells.Item(6, 11).Value = [double]"-38.02402990805938"
$ws.Cells.Item(6, 12).Value = [double]"12.80821536923895"
$ws.Cells.Item(6, 13).Value = [double]"0.001267901678596098"
$ws.Cells.Item(7, 2).Value = [double]"0.03644711399144783"
$ws.Cells.Item(7, 3).Value = [double]"-0.01671647636692574"
$ws.Cells.Item(7, 4).Value = '+'
$ws.Cells.Item(7, 5).Value = [double]"0.09797251220250203"
$ws.Cells.Item(7, 7).ClearContents()
$ws.Cells.Item(7, 8).Value = [double]"0.6249426245980445"
$ws.Cells.Item(7, 10).Value = [double]"26.61337885515504"
$ws.Cells.Item(7, 14).Value = [double]"0.5713629995406223"
$ws.Cells.Item(8, 2).Value = [double]"0.01058054056923997"
$ws.Cells.Item(8, 3).Value = [double]"-0.1386414594274318"
$ws.Cells.Item(8, 4).Value = '+'
$ws.Cells.Item(8, 6).ClearContents()
$ws.Cells.Item(8, 7).Value = '+'
$ws.Cells.Item(8, 8).Value = [double]"0.6180524392233377"
$ws.Cells.Item(8, 10).Value = [double]"26.31300776670598"
$ws.Cells.Item(8, 11).Value = [double]"-37.39524630264273"
$ws.Cells.Item(8, 12).Value = [double]"13.43699897465559"
$ws.Cells.Item(8, 13).Value = [double]"0.0009258634751589842"
$ws.Cells.Item(9, 2).Value = [double]"-0.03345313414186778"
$ws.Cells.Item(9, 3).Value = [double]"-0.0727078694538441"
$ws.Cells.Item(9, 4).Value = '+'
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(9, 6).Value = '+'
$ws.Cells.Item(9, 10).Value = [double]"24.91332396722979"
$ws.Cells.Item(9, 11).Value = [double]"-34.59587870369035"
$ws.Cells.Item(9, 12).Value = [double]"16.23636657360798"
$ws.Cells.Item(9, 13).Value = [double]"0.0002283873267747919"
$ws.Cells.Item(10, 2).ClearContents()
$ws.Cells.Item(10, 3).Value = [double]"-0.05520576681292619"
$ws.Cells.Item(10, 4).Value = '+'
$ws.Cells.Item(10, 10).Value = [double]"20.86002074549945"
$ws.Cells.Item(10, 11).Value = [double]"-32.29147006242748"
$ws.Cells.Item(10, 12).Value = [double]"18.54077521487085"
$ws.Cells.Item(10, 13).Value = [double]"7.215659361470195e-05"
$ws.Cells.Item(11, 2).Value = [double]"0.03373103834351063"
$ws.Cells.Item(11, 3).Value = [double]"-0.05373358937957363"
$ws.Cells.Item(11, 4).Value = '+'
$ws.Cells.Item(11, 11).Value = [double]"-31.59815783865717"
$ws.Cells.Item(11, 12).Value = [double]"19.23408743864115"
$ws.Cells.Item(11, 13).Value = [double]"5.10182063742772e-05"
$ws.Cells.Item(12, 2).ClearContents()
$ws.Cells.Item(12, 4).Value = '+'
$ws.Cells.Item(12, 8).Value = [double]"0.3796497006378914"
$ws.Cells.Item(12, 10).Value = [double]"18.31049137018108"
$ws.Cells.Item(12, 11).Value = [double]"-29.7933965334656"
$ws.Cells.Item(12, 12).Value = [double]"21.03884874383273"
$ws.Cells.Item(12, 13).Value = [double]"2.069313295967353e-05"
$ws.Cells.Item(12, 14).Value = [double]"0.3596384006584685"
$ws.Cells.Item(13, 2).Value = [double]"0.03609894030286132"
$ws.Cells.Item(13, 3).ClearContents()
$ws.Cells.Item(13, 4).Value = '+'
$ws.Cells.Item(13, 8).Value = [double]"0.4172615507292327"
$ws.Cells.Item(13, 10).Value = [double]"19.34249807354961"
$ws.Cells.Item(13, 11).Value = [double]"-29.25642471852779"
$ws.Cells.Item(13, 12).Value = [double]"21.57582055877054"
$ws.Cells.Item(13, 13).Value = [double]"1.582064917667436e-05"
$ws.Cells.Item(13, 14).Value = [double]"0.3784123207778481"
$ws.Cells.Item(14, 2).Value = [double]"-0.001765138510898472"
$ws.Cells.Item(14, 3).ClearContents()
$ws.Cells.Item(14, 4).Value = '+'
$ws.Cells.Item(14, 5).ClearContents()
$ws.Cells.Item(14, 6).Value = '+'
$ws.Cells.Item(14, 8).Value = [double]"0.4477026615080643"
$ws.Cells.Item(14, 10).Value = [double]"20.22775440220348"
$ws.Cells.Item(14, 11).Value = [double]"-28.23328658218474"
$ws.Cells.Item(14, 12).Value = [double]"22.59895869511359"
$ws.Cells.Item(14, 13).Value = [double]"9.485335064748589e-06"
$ws.Cells.Item(14, 14).Value = [double]"0.390568454077864"
$ws.Cells.Item(15, 3).Value = [double]"-0.06090616442362588"
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(15, 8).Value = [double]"0.1083900374261652"
$ws.Cells.Item(15, 10).Value = [double]"12.32520781637782"
$ws.Cells.Item(15, 11).Value = [double]"-17.82282942585909"
$ws.Cells.Item(15, 12).Value = [double]"33.00941585143924"
$ws.Cells.Item(15, 13).Value = [double]"5.205358294719213e-08"
$ws.Cells.Item(15, 14).Value = [double]"0.07962842573023499"
$ws.Cells.Item(16, 2).Value = [double]"0.04536254508886486"
$ws.Cells.Item(16, 3).Value = [double]"-0.05867458996061583"
$ws.Cells.Item(16, 4).ClearContents()
$ws.Cells.Item(16, 8).Value = [double]"0.1683704424322527"
$ws.Cells.Item(16, 13).Value = [double]"4.474003271893749e-08"
$ws.Cells.Item(16, 14).Value = [double]"0.1129284719277363"
$ws.Cells.Item(17, 10).Value = [double]"10.43222050504034"
$ws.Cells.Item(17, 11).Value = [double]"-16.46444101008069"
$ws.Cells.Item(17, 12).Value = [double]"34.36780426721764"
$ws.Cells.Item(17, 13).Value = [double]"2.639248792708289e-08"
$ws.Cells.Item(18, 2).Value = [double]"0.04824899505165226"
$ws.Cells.Item(18, 3).ClearContents()
$ws.Cells.Item(18, 8).Value = [double]"0.06802106217284964"
$ws.Cells.Item(18, 10).Value = [double]"11.59456405196224"
$ws.Cells.Item(18, 12).Value = [double]"34.47070338027039"
$ws.Cells.Item(18, 13).Value = [double]"2.506894589754373e-08"
$ws.Cells.Item(18, 14).Value = [double]"0.03795722546874802"
$ws.Cells.Item(19, 2).Value = [double]"0.04632972952549103"
$ws.Cells.Item(19, 3).Value = [double]"-0.0531323373367344"
$ws.Cells.Item(19, 4).ClearContents()
$ws.Cells.Item(19, 5).Value = [double]"0.01527973073000808"
$ws.Cells.Item(19, 7).ClearContents()
$ws.Cells.Item(19, 8).Value = [double]"0.1721062641606441"
$ws.Cells.Item(19, 10).Value = [double]"13.54858327969539"
$ws.Cells.Item(19, 11).Value = [double]"-14.87494433716856"
$ws.Cells.Item(19, 12).Value = [double]"35.95730094012977"
$ws.Cells.Item(19, 13).Value = [double]"1.192135208302038e-08"
$ws.Cells.Item(19, 14).Value = [double]"0.08646208459105553"

Write-Output "applied 152 changes"